# "more fron tend design" -- wipe the sample transaction rows back to a
# blank template and stretch the sheet down to 9 rows so there's extra
# blank room underneath for the next design pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("d"/"d" sample customer) - clear name/email/quantities/total.
$ws.Range("B2:C3").ClearContents()
$ws.Range("D2:H3").ClearContents()

# Row 4's leftover order number.
$ws.Range("A4").ClearContents()

# Stretch the sheet down to row 9 with blank rows so there's room to keep
# designing below the existing data.
$ws.Rows("7:7").Hidden = $true
$ws.Rows("7:7").Hidden = $false
$ws.Rows("8:8").Hidden = $true
$ws.Rows("8:8").Hidden = $false
$ws.Range("A9").Borders.LineStyle = -4142
